# Update the "Förändrad" date column (column C) for all data rows
# from 2023-09-06 (OADate 45175) to 2023-09-08 (OADate 45177).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $v = $cell.Value()
    if ($v -ne $null) {
        if ($v -is [DateTime]) {
            $oa = $v.ToOADate()
        } else {
            $oa = $v
        }
        if ($oa -eq 45175) {
            $cell.Value = 45177
        }
    }
}
